# Generate Report for Handoff
# Updates the localization-status report after a new handoff was generated
# for the "low" priority files (356ed863, 5560ddf0, 56f3be6a, fe708d45),
# moving them into the "ht" (handoff triggered) state with a fresh
# Latest Handoff Datetime, for both the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 correspond to files that just had a new handoff
# generated. Priority column (E) flips from "low" to "ht", and the
# Latest Handoff Datetime column (H) is refreshed.
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-30 16:38:20"
}

# de-de sheet: same set of rows / files, different handoff timestamp.
foreach ($r in 4..7) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-30 16:38:26"
}

# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# most recent per-language handoff datetime (de-de, in this dataset), so it
# advances to the same new timestamp as de-de's refreshed handoff date.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-08-30 16:38:26"
}
